# 05-table-features.xlsx: add a third "AutoWidthFalse" sheet/table-features
# fixture that exercises the corrected (non-auto-width) default column
# behaviour. Mirrors the existing "AutoWidthAll"/"AutoWidthBody" sheets but
# is appended after them, keeps Excel's default column widths (i.e. no
# explicit <cols> overrides), and carries its own short header/body text.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "AutoWidthFalse"

$ws.Range("A1").Value = "短"
$ws.Range("B1").Value = "長いヘッダーテキスト"
$ws.Range("A2").Value = "A"
$ws.Range("B2").Value = "データ"
